# REVER_DailyTracker_MONISHA.xlsx - "Add files via upload"
# Target sheet: NOV-2020
#   - D18 text changes from "Sick Leave" to "Leave"
#   - Two new rows (20, 21) are appended for Nov 19 & Nov 20, 2020
#   - Sheet scroll/selection moves down to the new bottom of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")
$ws.Activate()

# --- Update existing cell: "Sick Leave" -> "Leave" ---
# (Do this first so the shared-string table gets "Leave" allocated before
#  the two brand-new strings used below, matching the original edit order.)
$ws.Range("D18").Value = "Leave"

# --- Append row 20: Thu 19-Nov-2020 ---
$ws.Range("A19:G19").Copy()
$ws.Range("A20:G20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A21:G21").PasteSpecial(-4122)   # xlPasteFormats (row 21 too)
$excel.CutCopyMode = $false

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 44154
$ws.Range("C20").Value = "nMVAR "
$ws.Range("D20").Value = "QA_import_bat , QA_export_bat - tested"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "COMPLETED"

# --- Append row 21: Fri 20-Nov-2020 ---
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 44155
$ws.Range("C21").Value = "nMVAR "
$ws.Range("D21").Value = "nMVAR_QA tested fully"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "COMPLETED"

# --- Move viewport / selection to the new bottom of the sheet ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D26").Select() | Out-Null
